$wb = $excel.ActiveWorkbook

# SchemaOrganization
$ws = $wb.Worksheets.Item("SchemaOrganization")
$ws.Range("B3").Value = "http://example.com/organization3:Image0"

# RightsStatementsDotOrgRightsStatement
$ws = $wb.Worksheets.Item("RightsStatementsDotOrgRightsStatement")
$ws.Range("E2").Value = "Unless expressly stated otherwise, the organization that has made this Item available makes no warranties about the Item and cannot guarantee the accuracy of this Rights Statement. You are responsible for your own use."

# FoafPerson
$ws = $wb.Worksheets.Item("FoafPerson")
$ws.Range("E4").Value = "http://example.com/person4:Image1"

# SchemaPerson
$ws = $wb.Worksheets.Item("SchemaPerson")
$ws.Range("D2").Value = "http://example.com/person1:Image0"
$ws.Range("D3").Value = "http://example.com/person3:Image0"
$ws.Range("D4").Value = "http://example.com/person5:Image1"

# RdfProperty
$ws = $wb.Worksheets.Item("RdfProperty")
$ws.Range("C3").Value = "dcterms:extent:Image0"
$ws.Range("C7").Value = "dcterms:source:Image0"
$ws.Range("C10").Value = "dcterms:title:Image0"

# SkosConcept
$ws = $wb.Worksheets.Item("SkosConcept")
$ws.Range("B5").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:3:Image0"
$ws.Range("B7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:5:Image1"
$ws.Range("B10").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:8:Image0"
$ws.Range("B11").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:9:Image1"
$ws.Range("B12").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:10:Image1"
$ws.Range("B13").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:11:Image1"
$ws.Range("B14").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:12:Image1"
$ws.Range("B18").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:16:Image0"
$ws.Range("B19").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:17:Image0"
$ws.Range("B20").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:18:Image1"
$ws.Range("B22").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:20:Image1"
$ws.Range("B24").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:22:Image1"
$ws.Range("B27").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:25:Image0"
$ws.Range("B29").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:27:Image1"
$ws.Range("B31").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:29:Image0"
$ws.Range("B32").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:30:Image1"
$ws.Range("B33").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:31:Image1"
$ws.Range("B34").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:32:Image1"
$ws.Range("B35").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:33:Image0"
$ws.Range("B39").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:37:Image0"
$ws.Range("B43").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:41:Image1"
$ws.Range("B48").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:46:Image1"
$ws.Range("B49").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:47:Image0"
$ws.Range("B50").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:48:Image1"
$ws.Range("B51").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:49:Image1"
$ws.Range("B52").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:50:Image0"
$ws.Range("B54").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:52:Image1"
$ws.Range("B58").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:56:Image1"
$ws.Range("B59").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:57:Image0"
$ws.Range("B61").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:59:Image0"
$ws.Range("B62").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:60:Image1"
$ws.Range("B64").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:62:Image0"
$ws.Range("B68").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:66:Image0"
$ws.Range("B70").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:68:Image1"
$ws.Range("B73").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:71:Image0"
$ws.Range("B74").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:72:Image1"
$ws.Range("B76").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:74:Image0"
$ws.Range("B77").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:75:Image0"
$ws.Range("B78").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:76:Image1"

# SchemaDefinedTerm
$ws = $wb.Worksheets.Item("SchemaDefinedTerm")
$ws.Range("B2").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:80:Image0"
$ws.Range("B7").Value = "urn:paradicms_etl:pipeline:synthetic_data:concept:85:Image0"

# FoafOrganization
$ws = $wb.Worksheets.Item("FoafOrganization")
$ws.Range("C2").Value = "http://example.com/organization0:Image1"

# CreativeCommonsLicense: swap row 3 (Public Domain Mark 1.0) and row 4 (CC BY-SA 2.0) contents
$ws = $wb.Worksheets.Item("CreativeCommonsLicense")
$ws.Range("A3").Value = "http://creativecommons.org/licenses/by-sa/2.0/"
$ws.Range("E3").ClearContents()
$ws.Range("F3").ClearContents()
$ws.Range("G3").Value = "by-sa"
$ws.Range("H3").Value = "Attribution-ShareAlike 2.0 Generic"
$ws.Range("I3").Value = "'2.0"
$ws.Range("J3").NumberFormat = "General"
$ws.Range("J3").Value = ""
$ws.Range("K3").NumberFormat = "General"
$ws.Range("K3").Value = ""

$ws.Range("A4").Value = "http://creativecommons.org/publicdomain/mark/1.0/"
$ws.Range("E4").Value = "mark"
$ws.Range("F4").Value = "Public Domain Mark 1.0"
$ws.Range("G4").Value = "'1.0"
$ws.Range("H4").ClearContents()
$ws.Range("I4").ClearContents()
$ws.Range("J4").ClearContents()
$ws.Range("K4").ClearContents()

